$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 172-173; this shifts the existing rows 172-196 down to 174-198.
$ws.Rows("172:173").Insert()

# Fill in row 172 (new record)
$ws.Range("A172").Value = 6
$ws.Range("B172").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C172").Value = "Metropolitana"
$ws.Range("D172").Value = 44504
$ws.Range("E172").Value = 13
$ws.Range("F172").Value = "Fruta"
$ws.Range("G172").Value = 100101
$ws.Range("H172").Value = "Berries"
$ws.Range("I172").Value = 100101001
$ws.Range("J172").Value = "Arándano (blue)"
$ws.Range("K172").Value = "Sin especificar"
$ws.Range("L172").Value = "Primera"
$ws.Range("M172").Value = 450
$ws.Range("N172").Value = 7000
$ws.Range("O172").Value = 7000
$ws.Range("P172").Value = 7000
$ws.Range("Q172").Value = "$/bandeja 2 kilos"
$ws.Range("R172").Value = "Provincia de Curicó"
$ws.Range("S172").Value = 3500
$ws.Range("T172").Value = 2

# Fill in row 173 (new record)
$ws.Range("A173").Value = 6
$ws.Range("B173").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C173").Value = "Metropolitana"
$ws.Range("D173").Value = 44504
$ws.Range("E173").Value = 13
$ws.Range("F173").Value = "Fruta"
$ws.Range("G173").Value = 100101
$ws.Range("H173").Value = "Berries"
$ws.Range("I173").Value = 100101001
$ws.Range("J173").Value = "Arándano (blue)"
$ws.Range("K173").Value = "Sin especificar"
$ws.Range("L173").Value = "Segunda"
$ws.Range("M173").Value = 150
$ws.Range("N173").Value = 6000
$ws.Range("O173").Value = 6000
$ws.Range("P173").Value = 6000
$ws.Range("Q173").Value = "$/bandeja 2 kilos"
$ws.Range("R173").Value = "Provincia de Curicó"
$ws.Range("S173").Value = 3000
$ws.Range("T173").Value = 2
